# Changing label names for units
# - adds unit qualifiers to a handful of parameter labels
# - removes the now-redundant "IK_pump_hs" (additional loading power
#   investment) column from the Heat Storage sheet
# - switches the active sheet from "Techologies" to "Heat Storage"

$wb = $excel.ActiveWorkbook

$wsPrices = $wb.Worksheets.Item("prices and emmision factors")
$wsParams = $wb.Worksheets.Item("financal and other parameteres")
$wsStorage = $wb.Worksheets.Item("Heat Storage")

# --- "financal and other parameteres" sheet -------------------------------
# Renewable factor label switches from "Total" to "Minimum".
$wsParams.Range("D2").Value = "Minimum Renewable Factor [0-1]"

# --- "prices and emmision factors" sheet ---------------------------------
# Label for the emission factor row gains an explicit unit.
$wsPrices.Range("C2").Value = "emission factor [tCO2/MWh]"
[void]$wsPrices.Range("G22").Select()

# --- "financal and other parameteres" sheet (continued) -------------------
# CO2 price label gains an explicit unit.
$wsParams.Range("A2").Value = "CO2 Price [EUR/tC02]"
[void]$wsParams.Range("D20").Select()

# --- "Heat Storage" sheet --------------------------------------------------
# Column H ("IK_pump_hs" - investment costs for additional loading power)
# is dropped entirely; the columns to its right (OP_fix_hs, LT_hs) shift
# left to fill the gap.
[void]$wsStorage.Columns.Item(8).Delete()

# Remaining unit labels get more precise units.
$wsStorage.Range("H2").Value = "OPEX fix [€/MWha]"
$wsStorage.Range("G2").Value = "Invesment costs for additional storage capacity  [€/MWh]"

# Data update for the unloading/loading power capacity values.
$wsStorage.Range("C3").Value = 80
$wsStorage.Range("D3").Value = 80

# Heat Storage becomes the active sheet (was Techologies), with a new
# selection.
$wsStorage.Activate()
[void]$wsStorage.Range("G13").Select()
